$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 624, shifting the existing rows (old 624..665) down to
# (625..666). This mirrors the daily auto-push log gaining one more entry
# for 2026/01/14 just before the pre-existing "future schedule" block.
$ws.Rows(624).Insert()

# Column A holds a text date like "2026/12/29"; assigning a date-shaped
# string to .Value auto-coerces it to a real Excel date serial, so force
# text format first, then reset the cell style back to Normal afterwards
# so no stray style index is left on the cell (matches the plain,
# style-less date cells used throughout the rest of the column).
$ws.Range("A624").NumberFormat = "@"
$ws.Range("A624").Value = "2026/01/14"
$ws.Range("A624").Style = "Normal"

$ws.Range("B624").Value = "水"
$ws.Range("C624").Value = 19
$ws.Range("D624").Value = 28
